# Updated symbol list on Sat Dec 17 13:27:30 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" column (D) on the active sheet for a batch of
# coins with newly-scraped values. Column D holds text-typed numeric
# strings (t="inlineStr"/"s" in the sheet, not numeric cells), so each
# new value is entered with a leading apostrophe to force Excel to keep
# it as text instead of silently re-typing it as a Number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = [ordered]@{
    2  = "237.92"
    3  = "21.77"
    4  = "5.446"
    6  = "6.482"
    8  = "0.7963"
    9  = "1.014"
    10 = "0.1397"
    11 = "0.07320"
    12 = "0.03178"
    13 = "0.02976"
    14 = "0.09245"
    15 = "0.001661"
    16 = "3.256"
    17 = "0.04762"
    18 = "0.0005743"
    19 = "0.006227"
    20 = "0.005091"
    22 = "0.0001502"
    24 = "3.919"
    25 = "2.203"
    40 = "0.04109"
    41 = "0.006964"
    42 = "0.1042"
    43 = "0.003014"
    44 = "0.008625"
    45 = "0.00005435"
    46 = "0.00000000750"
    48 = "0.03564"
    50 = "0.01010"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}
